# This workbook contains weekly price-report rows (2 rows per reporting
# date: "Primera" / "Segunda" quality). A new weekly entry was inserted
# right before the existing row 43, pushing all rows from the old row 43
# onward down by two rows (old row 43 -> new row 45, ..., old row 164 ->
# new row 166). The two newly inserted rows (43 and 44) are populated
# with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 43; this shifts every row that was at
# 43 or below down by two (43->45, 44->46, ..., 164->166) and keeps all
# of their existing values/formatting intact.
$ws.Rows("43:44").Insert()

# New row 43 ("Primera" quality) for the new reporting date.
$ws.Cells.Item(43,1).Value  = 1
$ws.Cells.Item(43,2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(43,3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(43,4).Value  = 44459
$ws.Cells.Item(43,5).Value  = 15
$ws.Cells.Item(43,6).Value  = 100114014
$ws.Cells.Item(43,7).Value  = 'Betarraga'
$ws.Cells.Item(43,8).Value  = 'Sin especificar'
$ws.Cells.Item(43,9).Value  = 'Primera'
$ws.Cells.Item(43,10).Value = 1200
$ws.Cells.Item(43,11).Value = 500
$ws.Cells.Item(43,12).Value = 550
$ws.Cells.Item(43,13).Value = 525
$ws.Cells.Item(43,14).Value = '$/paquete 4 unidades'
$ws.Cells.Item(43,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(43,16).Value = 131
$ws.Cells.Item(43,17).Value = 4
$ws.Cells.Item(43,18).Value = 'Hortaliza'

# New row 44 ("Segunda" quality) for the same new reporting date.
$ws.Cells.Item(44,1).Value  = 1
$ws.Cells.Item(44,2).Value  = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(44,3).Value  = 'Arica y Parinacota'
$ws.Cells.Item(44,4).Value  = 44459
$ws.Cells.Item(44,5).Value  = 15
$ws.Cells.Item(44,6).Value  = 100114014
$ws.Cells.Item(44,7).Value  = 'Betarraga'
$ws.Cells.Item(44,8).Value  = 'Sin especificar'
$ws.Cells.Item(44,9).Value  = 'Segunda'
$ws.Cells.Item(44,10).Value = 1200
$ws.Cells.Item(44,11).Value = 500
$ws.Cells.Item(44,12).Value = 550
$ws.Cells.Item(44,13).Value = 525
$ws.Cells.Item(44,14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(44,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(44,16).Value = 105
$ws.Cells.Item(44,17).Value = 5
$ws.Cells.Item(44,18).Value = 'Hortaliza'
